$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3695.8572
$ws.Range("I18").Value = 3695.8572
$ws.Range("K18").Value = 3695.8572
$ws.Range("M18").Value = -3411.8572

$ws.Range("H51").Value = 3005.4602
$ws.Range("I51").Value = 2948.7437
$ws.Range("J51").Value = 3097.625
$ws.Range("K51").Value = 2948.7437
$ws.Range("L51").Value = 3097.625
$ws.Range("M51").Value = -2464.7437
$ws.Range("N51").Value = -4065.625

$ws.Range("H62").Value = 3387.5
$ws.Range("I62").Value = 2550
$ws.Range("J62").Value = 3666.6667
$ws.Range("K62").Value = 2550
$ws.Range("L62").Value = 3666.6667
$ws.Range("M62").Value = -1926
$ws.Range("N62").Value = -4914.6667

$ws.Range("H65").Value = 3387.5
$ws.Range("I65").Value = 2550
$ws.Range("J65").Value = 3666.6667
$ws.Range("K65").Value = 12750
$ws.Range("L65").Value = 18333.3335
$ws.Range("M65").Value = -9630
$ws.Range("N65").Value = -24573.3335

$ws.Range("H98").Value = 691.24
$ws.Range("I98").Value = 644
$ws.Range("J98").Value = 3006
$ws.Range("K98").Value = 644
$ws.Range("L98").Value = 3006
$ws.Range("M98").Value = 854
$ws.Range("N98").Value = -6002

$ws.Range("H100").Value = 1859.64
$ws.Range("I100").Value = 1579.2667
$ws.Range("J100").Value = 2280.2
$ws.Range("K100").Value = 1579.2667
$ws.Range("L100").Value = 2280.2
$ws.Range("M100").Value = -1038.2667
$ws.Range("N100").Value = -3362.2

$ws.Range("H103").Value = 587.125
$ws.Range("I103").Value = 493.14285
$ws.Range("J103").Value = 1245
$ws.Range("K103").Value = 1479.42855
$ws.Range("L103").Value = 3735
$ws.Range("M103").Value = -893.4285500000001
$ws.Range("N103").Value = -4907

$ws.Range("H107").Value = 3723.4138
$ws.Range("I107").Value = 3563.0454
$ws.Range("K107").Value = 3563.0454
$ws.Range("M107").Value = -1643.0454

$ws.Range("H122").Value = 691.24
$ws.Range("I122").Value = 644
$ws.Range("J122").Value = 3006
$ws.Range("K122").Value = 1932
$ws.Range("L122").Value = 9018
$ws.Range("M122").Value = 518
$ws.Range("N122").Value = -13918

$ws.Range("H129").Value = 20908.2
$ws.Range("I129").Value = 13010.25
$ws.Range("K129").Value = 39030.75
$ws.Range("M129").Value = -34030.75

$ws.Range("H137").Value = 33345300
$ws.Range("I137").Value = 50002320
$ws.Range("K137").Value = 150006960
$ws.Range("M137").Value = -150004410

$ws.Range("H138").Value = 3561.6965
$ws.Range("I138").Value = 3766.6572
$ws.Range("J138").Value = 3220.0952
$ws.Range("K138").Value = 11299.9716
$ws.Range("L138").Value = 9660.285600000001
$ws.Range("M138").Value = -6159.971600000001
$ws.Range("N138").Value = -19940.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 609386.4
$ws.Range("I32").Value = 649760.6
$ws.Range("J32").Value = 10502.333
$ws.Range("K32").Value = 649760.6
$ws.Range("L32").Value = 10502.333
$ws.Range("M32").Value = -649473.6
$ws.Range("N32").Value = -11076.333

$ws.Range("H61").Value = 2636256
$ws.Range("I61").Value = 4786.5884
$ws.Range("J61").Value = 25003746
$ws.Range("K61").Value = 4786.5884
$ws.Range("L61").Value = 25003746
$ws.Range("M61").Value = -4574.5884
$ws.Range("N61").Value = -25004170

$ws.Range("H74").Value = 1137543
$ws.Range("I74").Value = 1793485.9
$ws.Range("K74").Value = 1793485.9
$ws.Range("M74").Value = -1792611.9

$ws.Range("H77").Value = 1137543
$ws.Range("I77").Value = 1793485.9
$ws.Range("K77").Value = 8967429.5
$ws.Range("M77").Value = -8963061.5

$ws.Range("H132").Value = 737199.2
$ws.Range("I132").Value = 834894.4
$ws.Range("J132").Value = 4485.25
$ws.Range("K132").Value = 2504683.2
$ws.Range("L132").Value = 13455.75
$ws.Range("M132").Value = -2502153.2
$ws.Range("N132").Value = -18515.75

$ws.Range("H136").Value = 2636256
$ws.Range("I136").Value = 4786.5884
$ws.Range("J136").Value = 25003746
$ws.Range("K136").Value = 14359.7652
$ws.Range("L136").Value = 75011238
$ws.Range("M136").Value = -11809.7652
$ws.Range("N136").Value = -75016338

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1404
$ws.Range("I107").Value = 1348.9656
$ws.Range("K107").Value = 1348.9656
$ws.Range("M107").Value = 571.0344

$ws.Range("H134").Value = 3405148.8
$ws.Range("I134").Value = 1790.0454
$ws.Range("K134").Value = 5370.1362
$ws.Range("M134").Value = -2835.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2527271.2
$ws.Range("I31").Value = 2779814
$ws.Range("J31").Value = 1842
$ws.Range("K31").Value = 2779814
$ws.Range("L31").Value = 1842
$ws.Range("M31").Value = -2779519
$ws.Range("N31").Value = -2432

$ws.Range("H34").Value = 2527271.2
$ws.Range("I34").Value = 2779814
$ws.Range("J34").Value = 1842
$ws.Range("K34").Value = 2779814
$ws.Range("L34").Value = 1842
$ws.Range("M34").Value = -2779612
$ws.Range("N34").Value = -2246

$ws.Range("H58").Value = 1815690.5
$ws.Range("I58").Value = 1962.1875
$ws.Range("J58").Value = 5961355
$ws.Range("K58").Value = 1962.1875
$ws.Range("L58").Value = 5961355
$ws.Range("M58").Value = -1759.1875
$ws.Range("N58").Value = -5961761

$ws.Range("H122").Value = 8144.439
$ws.Range("I122").Value = 2033
$ws.Range("J122").Value = 19930.785
$ws.Range("K122").Value = 6099
$ws.Range("L122").Value = 59792.355
$ws.Range("M122").Value = -3649
$ws.Range("N122").Value = -64692.355

$ws.Range("H132").Value = 1684.4706
$ws.Range("I132").Value = 1494.9546
$ws.Range("J132").Value = 2875.7144
$ws.Range("K132").Value = 4484.8638
$ws.Range("L132").Value = 8627.143199999999
$ws.Range("M132").Value = -1954.8638
$ws.Range("N132").Value = -13687.1432

$ws.Range("H134").Value = 1178.1765
$ws.Range("I134").Value = 1062.3636
$ws.Range("K134").Value = 3187.0908
$ws.Range("M134").Value = -652.0907999999999

$ws.Range("H136").Value = 1815690.5
$ws.Range("I136").Value = 1962.1875
$ws.Range("J136").Value = 5961355
$ws.Range("K136").Value = 5886.5625
$ws.Range("L136").Value = 17884065
$ws.Range("M136").Value = -3336.5625
$ws.Range("N136").Value = -17889165

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2127158.5
$ws.Range("I4").Value = 2714229.5
$ws.Range("K4").Value = 8142688.5
$ws.Range("M4").Value = -8142576.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7433.442
$ws.Range("I132").Value = 7729.025
$ws.Range("K132").Value = 23187.075
$ws.Range("M132").Value = -20657.075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1906419.4
$ws.Range("I132").Value = 3031811.8
$ws.Range("J132").Value = 1909.0769
$ws.Range("K132").Value = 9095435.399999999
$ws.Range("L132").Value = 5727.2307
$ws.Range("M132").Value = -9092905.399999999
$ws.Range("N132").Value = -10787.2307

$ws.Range("H140").Value = 50000000
$ws.Range("J140").Value = 50000000
$ws.Range("L140").Value = 50000000
$ws.Range("N140").Value = -50010360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3788885.2
$ws.Range("I132").Value = 3876996.5
$ws.Range("J132").Value = 99
$ws.Range("K132").Value = 11630989.5
$ws.Range("L132").Value = 297
$ws.Range("M132").Value = -11628459.5
$ws.Range("N132").Value = -5357

$ws.Range("H136").Value = 34541924
$ws.Range("I136").Value = 34783520
$ws.Range("J136").Value = 33333954
$ws.Range("K136").Value = 104350560
$ws.Range("L136").Value = 100001862
$ws.Range("M136").Value = -104348010
$ws.Range("N136").Value = -100006962
